# Refresh the crypto price table (cryptos.xlsx) with freshly scraped values.
# Sheet layout: A=rank index, B=Coin, C=Link, D=Price, E=Volume(1h).
# Every populated cell in the sheet was originally written as literal text
# (prices use "." as a thousands separator, e.g. "26.895.60"), so every
# replacement below must also land as text, never as a re-interpreted number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new price strings (e.g. "6.366") are valid Excel numbers, and a
# plain Value assignment would let Excel silently convert them into floats -
# dropping meaningful trailing zeros (e.g. "0.7470" -> 0.747). Force those
# specific cells to Text format first so the exact string is preserved, just
# like every other (non-numeric-looking) price cell already in the column.
$textFormatCells = @(
    "D5", "D7", "D8", "D9", "D10", "D14", "D15", "D17",
    "D18", "D24", "D25", "D27", "D28", "D29", "D30", "D32",
    "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41",
    "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D51"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.887.54"
$ws.Range("E2").Value = "  -0.93%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.868.03"
$ws.Range("E3").Value = "  +0.04%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "305.92"
$ws.Range("E5").Value = "  -0.03%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.03%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.5098"
$ws.Range("E7").Value = "  -1.03%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3657"
$ws.Range("E8").Value = "  -2.83%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.07169"
$ws.Range("E9").Value = "  +0.28%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "0.8897"
$ws.Range("E10").Value = "  -0.19%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.86%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.871.87"
$ws.Range("E13").Value = "  +0.24%  "

# Row 14 - Litecoin
$ws.Range("D14").Value = "94.52"
$ws.Range("E14").Value = "  +5.56%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "5.217"
$ws.Range("E15").Value = "  -1.73%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  -0.01%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.000008487"
$ws.Range("E17").Value = "  +0.15%  "

# Row 18 - Avalanche
$ws.Range("D18").Value = "14.13"
$ws.Range("E18").Value = "  +0.21%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.05%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "26.938.69"
$ws.Range("E20").Value = "  -0.88%  "

# Row 22 - Wrapped liquid staked Ether 2.0
$ws.Range("D22").Value = "2.110.17"
$ws.Range("E22").Value = "  +0.64%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  -1.35%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "6.366"

# Row 25 - Monero
$ws.Range("D25").Value = "147.64"
$ws.Range("E25").Value = "  +1.25%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  -3.12%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "17.85"
$ws.Range("E27").Value = "  -0.71%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "2.081"
$ws.Range("E28").Value = "  -0.17%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "113.37"
$ws.Range("E29").Value = "  +0.46%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "4.670"
$ws.Range("E30").Value = "  +0.23%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  +0.42%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "0.09128"
$ws.Range("E32").Value = "  -1.46%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.05024"
$ws.Range("E33").Value = "  -1.86%  "

# Row 34 / Row 35 - HuobiToken and ImmutableX swapped rank positions
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7470"
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.976"
$ws.Range("E35").Value = "  -3.42%  "

# Row 36 - ARBITRUM
$ws.Range("D36").Value = "1.150"
$ws.Range("E36").Value = "  -0.85%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "3.225"
$ws.Range("E37").Value = "  +3.97%  "

# Row 38 - RenderToken
$ws.Range("D38").Value = "2.511"
$ws.Range("E38").Value = "  +0.36%  "

# Row 39 - TheSandbox
$ws.Range("D39").Value = "0.5589"
$ws.Range("E39").Value = "  +5.77%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -2.27%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "1.071"
$ws.Range("E41").Value = "  -0.68%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "6.579"
$ws.Range("E42").Value = "  +1.07%  "

# Row 43 - Quant
$ws.Range("D43").Value = "115.44"
$ws.Range("E43").Value = "  -1.02%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "8.547"
$ws.Range("E44").Value = "  +2.72%  "

# Row 45 - Algorand
$ws.Range("E45").Value = "  +0.86%  "

# Row 46 - Decentraland
$ws.Range("D46").Value = "0.4764"
$ws.Range("E46").Value = "  +3.18%  "

# Row 47 - PaxDollar
$ws.Range("D47").Value = "0.9996"
$ws.Range("E47").Value = "  +0.04%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "10.06"
$ws.Range("E48").Value = "  +0.76%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "1.552"
$ws.Range("E49").Value = "  -0.48%  "

# Row 50 - Elrond
$ws.Range("E50").Value = "  +1.05%  "

# Row 51 - Aave
$ws.Range("D51").Value = "62.95"
$ws.Range("E51").Value = "  -1.13%  "
